# Add two new data rows (32 and 33) to the "Artfynd" sheet, matching the
# rows already present in the table (same columns / layout as row 31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 32 ----
$ws.Range("A32").Value = 112213272
$ws.Range("B32").Value = 89405
$ws.Range("C32").Value = "Ovaliderad"
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 1202
$ws.Range("F32").Value = "Ullticka"
$ws.Range("G32").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H32").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("P32").Value = "Simsbodarna O, Dlr"
$ws.Range("Q32").Value = 515738
$ws.Range("R32").Value = 6704726
$ws.Range("S32").Value = 1
$ws.Range("T32").Value = "Dalarna"
$ws.Range("U32").Value = "Borlänge"
$ws.Range("V32").Value = "Dalarna"
$ws.Range("W32").Value = "Stora Tuna"

# Date/time columns are stored as plain text in the workbook, not as real
# Excel date/time values, so force the cells to Text format before typing
# the value (otherwise Excel auto-converts "2023-09-20" into a date
# serial number), then drop the formatting again so no new style is left
# applied to the cell.
$ws.Range("Y32").NumberFormat = "@"
$ws.Range("AA32").NumberFormat = "@"
$ws.Range("Y32").Value = "2023-09-20"
$ws.Range("Z32").Value = "13:14"
$ws.Range("AA32").Value = "2023-09-20"
$ws.Range("AB32").Value = "13:14"
$ws.Range("Y32:AB32").ClearFormats()

$ws.Range("AD32").Value = $false
$ws.Range("AE32").Value = $false
$ws.Range("AG32").Value = $false
$ws.Range("AW32").Value = "Lars-Erik Nilsson"
$ws.Range("AX32").Value = "Lars-Erik Nilsson"

# ---- Row 33 ----
$ws.Range("A33").Value = 112213305
$ws.Range("B33").Value = 89369
$ws.Range("C33").Value = "Ovaliderad"
$ws.Range("D33").Value = "LC"
$ws.Range("E33").Value = 5447
$ws.Range("F33").Value = "Vedticka"
$ws.Range("G33").Value = "Fuscoporia viticola"
$ws.Range("H33").Value = "(Schwein.) Murrill"
$ws.Range("P33").Value = "Simsbodarna O, Dlr"
$ws.Range("Q33").Value = 515748
$ws.Range("R33").Value = 6704727
$ws.Range("S33").Value = 1
$ws.Range("T33").Value = "Dalarna"
$ws.Range("U33").Value = "Borlänge"
$ws.Range("V33").Value = "Dalarna"
$ws.Range("W33").Value = "Stora Tuna"

$ws.Range("Y33").NumberFormat = "@"
$ws.Range("AA33").NumberFormat = "@"
$ws.Range("Y33").Value = "2023-09-20"
$ws.Range("Z33").Value = "13:14"
$ws.Range("AA33").Value = "2023-09-20"
$ws.Range("AB33").Value = "13:14"
$ws.Range("Y33:AB33").ClearFormats()

$ws.Range("AD33").Value = $false
$ws.Range("AE33").Value = $false
$ws.Range("AG33").Value = $false
$ws.Range("AW33").Value = "Lars-Erik Nilsson"
$ws.Range("AX33").Value = "Lars-Erik Nilsson"
